# Auto-generated Excel COM-interop script
# Applies numeric corrections to currentAveragePrice / Leve price / profit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, as produced by the
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 10443.929
$ws.Range("I106").Value = 6162.8887
$ws.Range("J106").Value = 18149.8
$ws.Range("K106").Value = 6162.8887
$ws.Range("L106").Value = 18149.8
$ws.Range("M106").Value = -5531.8887
$ws.Range("N106").Value = -19411.8
$ws.Range("H132").Value = 4629.8423
$ws.Range("I132").Value = 4362.7417
$ws.Range("J132").Value = 5812.7144
$ws.Range("K132").Value = 13088.2251
$ws.Range("L132").Value = 17438.1432
$ws.Range("M132").Value = -10558.2251
$ws.Range("N132").Value = -22498.1432
$ws.Range("H137").Value = 3473.0505
$ws.Range("I137").Value = 1619
$ws.Range("J137").Value = 3942.4304
$ws.Range("K137").Value = 4857
$ws.Range("L137").Value = 11827.2912
$ws.Range("M137").Value = -2307
$ws.Range("N137").Value = -16927.2912
$ws.Range("H138").Value = 2893.4167
$ws.Range("I138").Value = 1796.5555
$ws.Range("J138").Value = 3551.5334
$ws.Range("K138").Value = 5389.666499999999
$ws.Range("L138").Value = 10654.6002
$ws.Range("M138").Value = -249.6664999999994
$ws.Range("N138").Value = -20934.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 218933.11
$ws.Range("I32").Value = 228679.75
$ws.Range("J32").Value = 4507
$ws.Range("K32").Value = 228679.75
$ws.Range("L32").Value = 4507
$ws.Range("M32").Value = -228392.75
$ws.Range("N32").Value = -5081
$ws.Range("H61").Value = 1930.5714
$ws.Range("I61").Value = 1993.5454
$ws.Range("J61").Value = 1699.6666
$ws.Range("K61").Value = 1993.5454
$ws.Range("L61").Value = 1699.6666
$ws.Range("M61").Value = -1781.5454
$ws.Range("N61").Value = -2123.6666
$ws.Range("H132").Value = 25642724
$ws.Range("I132").Value = 31251584
$ws.Range("J132").Value = 2217
$ws.Range("K132").Value = 93754752
$ws.Range("L132").Value = 6651
$ws.Range("M132").Value = -93752222
$ws.Range("N132").Value = -11711
$ws.Range("H136").Value = 1930.5714
$ws.Range("I136").Value = 1993.5454
$ws.Range("J136").Value = 1699.6666
$ws.Range("K136").Value = 5980.6362
$ws.Range("L136").Value = 5098.9998
$ws.Range("M136").Value = -3430.6362
$ws.Range("N136").Value = -10198.9998
$ws.Range("H141").Value = 20429
$ws.Range("J141").Value = 20429
$ws.Range("L141").Value = 20429
$ws.Range("N141").Value = -30789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 17145.666
$ws.Range("J64").Value = 25445.25
$ws.Range("L64").Value = 25445.25
$ws.Range("N64").Value = -25895.25
$ws.Range("H67").Value = 17145.666
$ws.Range("J67").Value = 25445.25
$ws.Range("L67").Value = 25445.25
$ws.Range("N67").Value = -27005.25
$ws.Range("H99").Value = 1198.5
$ws.Range("I99").Value = 1042
$ws.Range("J99").Value = 1355
$ws.Range("K99").Value = 1042
$ws.Range("L99").Value = 1355
$ws.Range("M99").Value = 456
$ws.Range("N99").Value = -4351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4423.5293
$ws.Range("I31").Value = 1787.9166
$ws.Range("J31").Value = 4988.3037
$ws.Range("K31").Value = 1787.9166
$ws.Range("L31").Value = 4988.3037
$ws.Range("M31").Value = -1492.9166
$ws.Range("N31").Value = -5578.3037
$ws.Range("H34").Value = 4423.5293
$ws.Range("I34").Value = 1787.9166
$ws.Range("J34").Value = 4988.3037
$ws.Range("K34").Value = 1787.9166
$ws.Range("L34").Value = 4988.3037
$ws.Range("M34").Value = -1585.9166
$ws.Range("N34").Value = -5392.3037
$ws.Range("H99").Value = 399640.16
$ws.Range("I99").Value = 927827.8
$ws.Range("J99").Value = 3499.375
$ws.Range("K99").Value = 927827.8
$ws.Range("L99").Value = 3499.375
$ws.Range("M99").Value = -926329.8
$ws.Range("N99").Value = -6495.375
$ws.Range("H122").Value = 4173.857
$ws.Range("I122").Value = 4536.1665
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 13608.4995
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -11158.4995
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 399640.16
$ws.Range("I126").Value = 927827.8
$ws.Range("J126").Value = 3499.375
$ws.Range("K126").Value = 2783483.4
$ws.Range("L126").Value = 10498.125
$ws.Range("M126").Value = -2781013.4
$ws.Range("N126").Value = -15438.125
$ws.Range("H132").Value = 956771
$ws.Range("I132").Value = 560010.75
$ws.Range("J132").Value = 3337332.8
$ws.Range("K132").Value = 1680032.25
$ws.Range("L132").Value = 10011998.4
$ws.Range("M132").Value = -1677502.25
$ws.Range("N132").Value = -10017058.4
$ws.Range("H134").Value = 2400.7097
$ws.Range("I134").Value = 1702.0476
$ws.Range("J134").Value = 3867.9
$ws.Range("K134").Value = 5106.142800000001
$ws.Range("L134").Value = 11603.7
$ws.Range("M134").Value = -2571.142800000001
$ws.Range("N134").Value = -16673.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78
$ws.Range("I12").Value = 58.833332
$ws.Range("J12").Value = 86.84614999999999
$ws.Range("K12").Value = 176.499996
$ws.Range("L12").Value = 260.53845
$ws.Range("M12").Value = -3.49999600000001
$ws.Range("N12").Value = -606.53845
$ws.Range("H68").Value = 2400
$ws.Range("I68").Value = 2033.3334
$ws.Range("J68").Value = 2840
$ws.Range("K68").Value = 6100.0002
$ws.Range("L68").Value = 8520
$ws.Range("M68").Value = -5289.0002
$ws.Range("N68").Value = -10142
$ws.Range("H71").Value = 2400
$ws.Range("I71").Value = 2033.3334
$ws.Range("J71").Value = 2840
$ws.Range("K71").Value = 18300.0006
$ws.Range("L71").Value = 25560
$ws.Range("M71").Value = -14244.0006
$ws.Range("N71").Value = -33672
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 253408.3
$ws.Range("I132").Value = 306009.25
$ws.Range("J132").Value = 5432.4287
$ws.Range("K132").Value = 918027.75
$ws.Range("L132").Value = 16297.2861
$ws.Range("M132").Value = -915497.75
$ws.Range("N132").Value = -21357.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10080.454
$ws.Range("I68").Value = 1950
$ws.Range("J68").Value = 11887.223
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 11887.223
$ws.Range("M68").Value = -1201
$ws.Range("N68").Value = -13385.223
$ws.Range("H71").Value = 10080.454
$ws.Range("I71").Value = 1950
$ws.Range("J71").Value = 11887.223
$ws.Range("K71").Value = 9750
$ws.Range("L71").Value = 59436.115
$ws.Range("M71").Value = -6006
$ws.Range("N71").Value = -66924.11499999999
$ws.Range("H132").Value = 7252.5
$ws.Range("I132").Value = 3718.75
$ws.Range("J132").Value = 10079.5
$ws.Range("K132").Value = 11156.25
$ws.Range("L132").Value = 30238.5
$ws.Range("M132").Value = -8626.25
$ws.Range("N132").Value = -35298.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 172511.5
$ws.Range("I45").Value = 69420
$ws.Range("J45").Value = 206875.33
$ws.Range("K45").Value = 69420
$ws.Range("L45").Value = 206875.33
$ws.Range("M45").Value = -68929
$ws.Range("N45").Value = -207857.33
$ws.Range("H132").Value = 5123903
$ws.Range("I132").Value = 5953914.5
$ws.Range("K132").Value = 17861743.5
$ws.Range("M132").Value = -17859213.5
$ws.Range("H133").Value = 130353
$ws.Range("J133").Value = 130353
$ws.Range("L133").Value = 130353
$ws.Range("N133").Value = -140473

